$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.312883734703064
$ws.Range("B1").Value = 4.11850118637085
$ws.Range("C1").Value = 6.002137184143066
$ws.Range("D1").Value = 1.481329202651978
$ws.Range("E1").Value = 0.8344982266426086
